$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.210.37"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.860.35"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2855"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07893"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "1.866.28"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.162"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6791"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "30.204.07"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +7.27%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.367"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "2.109.01"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007300"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.162"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.928"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.379"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09719"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.361"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.478"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04723"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7078"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.324"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8477"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "969.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.175"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.243"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05636"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
